# Auto-generated edit script applying numeric corrections to the Profits sheets.
# Each block targets one worksheet/row combination identified from the source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 5197
$ws.Range("J40").Value = 3998.8
$ws.Range("L40").Value = 3998.8
$ws.Range("N40").Value = -4348.8
# Row 64
$ws.Range("H64").Value = 5453.2856
$ws.Range("I64").Value = 6060.6665
$ws.Range("J64").Value = 4997.75
$ws.Range("K64").Value = 6060.6665
$ws.Range("L64").Value = 4997.75
$ws.Range("M64").Value = -5812.6665
$ws.Range("N64").Value = -5493.75
# Row 67
$ws.Range("H67").Value = 5453.2856
$ws.Range("I67").Value = 6060.6665
$ws.Range("J67").Value = 4997.75
$ws.Range("K67").Value = 6060.6665
$ws.Range("L67").Value = 4997.75
$ws.Range("M67").Value = -5202.6665
$ws.Range("N67").Value = -6713.75
# Row 86
$ws.Range("H86").Value = 24042790
$ws.Range("I86").Value = 12503298
$ws.Range("K86").Value = 12503298
$ws.Range("M86").Value = -12502175
# Row 89
$ws.Range("H89").Value = 24042790
$ws.Range("I89").Value = 12503298
$ws.Range("K89").Value = 62516490
$ws.Range("M89").Value = -62510874
# Row 137
$ws.Range("H137").Value = 3985.1035
$ws.Range("I137").Value = 1952
$ws.Range("J137").Value = 5055.1577
$ws.Range("K137").Value = 5856
$ws.Range("L137").Value = 15165.4731
$ws.Range("M137").Value = -3306
$ws.Range("N137").Value = -20265.4731
# Row 138
$ws.Range("H138").Value = 3507.5576
$ws.Range("J138").Value = 3983.9048
$ws.Range("L138").Value = 11951.7144
$ws.Range("N138").Value = -22231.7144

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 400.66666
$ws.Range("J5").Value = 400.66666
$ws.Range("L5").Value = 400.66666
$ws.Range("N5").Value = -624.66666
# Row 61
$ws.Range("H61").Value = 5186.5
$ws.Range("I61").Value = 4750
$ws.Range("K61").Value = 4750
$ws.Range("M61").Value = -4538
# Row 74
$ws.Range("H74").Value = 4453.0454
$ws.Range("I74").Value = 1870.4546
$ws.Range("K74").Value = 1870.4546
$ws.Range("M74").Value = -996.4546
# Row 77
$ws.Range("H77").Value = 4453.0454
$ws.Range("I77").Value = 1870.4546
$ws.Range("K77").Value = 9352.273000000001
$ws.Range("M77").Value = -4984.273000000001
# Row 86
$ws.Range("H86").Value = 150314
$ws.Range("J86").Value = 150314
$ws.Range("L86").Value = 150314
$ws.Range("N86").Value = -152686
# Row 89
$ws.Range("H89").Value = 150314
$ws.Range("J89").Value = 150314
$ws.Range("L89").Value = 450942
$ws.Range("N89").Value = -462798
# Row 132
$ws.Range("H132").Value = 2976.75
$ws.Range("I132").Value = 2484.5
$ws.Range("J132").Value = 4781.6665
$ws.Range("K132").Value = 7453.5
$ws.Range("L132").Value = 14344.9995
$ws.Range("M132").Value = -4923.5
$ws.Range("N132").Value = -19404.9995
# Row 136
$ws.Range("H136").Value = 5186.5
$ws.Range("I136").Value = 4750
$ws.Range("K136").Value = 14250
$ws.Range("M136").Value = -11700

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 400.66666
$ws.Range("J4").Value = 400.66666
$ws.Range("L4").Value = 400.66666
$ws.Range("N4").Value = -630.66666
# Row 86
$ws.Range("H86").Value = 2586.7856
$ws.Range("I86").Value = 2388.889
$ws.Range("J86").Value = 2943
$ws.Range("K86").Value = 2388.889
$ws.Range("L86").Value = 2943
$ws.Range("M86").Value = -1265.889
$ws.Range("N86").Value = -5189
# Row 89
$ws.Range("H89").Value = 2586.7856
$ws.Range("I89").Value = 2388.889
$ws.Range("J89").Value = 2943
$ws.Range("K89").Value = 11944.445
$ws.Range("L89").Value = 14715
$ws.Range("M89").Value = -6328.445
$ws.Range("N89").Value = -25947
# Row 97
$ws.Range("H97").Value = 8714.223
$ws.Range("I97").Value = 2759.7144
$ws.Range("J97").Value = 29555
$ws.Range("K97").Value = 2759.7144
$ws.Range("L97").Value = 29555
$ws.Range("M97").Value = -1768.7144
$ws.Range("N97").Value = -31537
# Row 107
$ws.Range("H107").Value = 3399.8948
$ws.Range("I107").Value = 3366.5
$ws.Range("J107").Value = 3493.4
$ws.Range("K107").Value = 3366.5
$ws.Range("L107").Value = 3493.4
$ws.Range("M107").Value = -1446.5
$ws.Range("N107").Value = -7333.4
# Row 134
$ws.Range("H134").Value = 14873.889
$ws.Range("I134").Value = 16441.066
$ws.Range("J134").Value = 7038
$ws.Range("K134").Value = 49323.198
$ws.Range("L134").Value = 21114
$ws.Range("M134").Value = -46788.198
$ws.Range("N134").Value = -26184

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 8547.312
$ws.Range("I31").Value = 1083.1794
$ws.Range("J31").Value = 21779.182
$ws.Range("K31").Value = 1083.1794
$ws.Range("L31").Value = 21779.182
$ws.Range("M31").Value = -788.1794
$ws.Range("N31").Value = -22369.182
# Row 34
$ws.Range("H34").Value = 8547.312
$ws.Range("I34").Value = 1083.1794
$ws.Range("J34").Value = 21779.182
$ws.Range("K34").Value = 1083.1794
$ws.Range("L34").Value = 21779.182
$ws.Range("M34").Value = -881.1794
$ws.Range("N34").Value = -22183.182
# Row 58
$ws.Range("H58").Value = 11314.814
$ws.Range("I58").Value = 22154.584
$ws.Range("J58").Value = 2643
$ws.Range("K58").Value = 22154.584
$ws.Range("L58").Value = 2643
$ws.Range("M58").Value = -21951.584
$ws.Range("N58").Value = -3049
# Row 122
$ws.Range("H122").Value = 6422.4287
$ws.Range("I122").Value = 2863.5
$ws.Range("J122").Value = 11167.667
$ws.Range("K122").Value = 8590.5
$ws.Range("L122").Value = 33503.001
$ws.Range("M122").Value = -6140.5
$ws.Range("N122").Value = -38403.001
# Row 132
$ws.Range("H132").Value = 1936.2727
$ws.Range("I132").Value = 2324.75
$ws.Range("J132").Value = 1714.2858
$ws.Range("K132").Value = 6974.25
$ws.Range("L132").Value = 5142.857400000001
$ws.Range("M132").Value = -4444.25
$ws.Range("N132").Value = -10202.8574
# Row 136
$ws.Range("H136").Value = 11314.814
$ws.Range("I136").Value = 22154.584
$ws.Range("J136").Value = 2643
$ws.Range("K136").Value = 66463.75199999999
$ws.Range("L136").Value = 7929
$ws.Range("M136").Value = -63913.75199999999
$ws.Range("N136").Value = -13029

$ws = $wb.Worksheets.Item("CUL")
# Row 70
$ws.Range("H70").Value = 3035.5454
$ws.Range("I70").Value = 350.25
$ws.Range("K70").Value = 1050.75
$ws.Range("M70").Value = -735.75
# Row 73
$ws.Range("H73").Value = 3035.5454
$ws.Range("I73").Value = 350.25
$ws.Range("K73").Value = 1050.75
$ws.Range("M73").Value = 41.25
# Row 93
$ws.Range("H93").Value = 5915.35
$ws.Range("J93").Value = 6016.1577
$ws.Range("L93").Value = 18048.4731
$ws.Range("N93").Value = -21792.4731
# Row 96
$ws.Range("H96").Value = 6747.5
$ws.Range("I96").Value = 4500
$ws.Range("K96").Value = 13500
$ws.Range("M96").Value = -11441
# Row 111
$ws.Range("H111").Value = 4384.6665
$ws.Range("I111").Value = 4384.6665
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 13153.9995
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -10086.9995
$ws.Range("N111").ClearContents()
# Row 137
$ws.Range("H137").Value = 3646.1333
$ws.Range("I137").Value = 1208.3334
$ws.Range("K137").Value = 3625.0002
$ws.Range("M137").Value = 1474.9998

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 110129.5
$ws.Range("I80").Value = 170592.33
$ws.Range("J80").Value = 49666.668
$ws.Range("K80").Value = 170592.33
$ws.Range("L80").Value = 49666.668
$ws.Range("M80").Value = -169594.33
$ws.Range("N80").Value = -51662.668
# Row 83
$ws.Range("H83").Value = 110129.5
$ws.Range("I83").Value = 170592.33
$ws.Range("J83").Value = 49666.668
$ws.Range("K83").Value = 852961.6499999999
$ws.Range("L83").Value = 248333.34
$ws.Range("M83").Value = -847969.6499999999
$ws.Range("N83").Value = -258317.34
# Row 122
$ws.Range("H122").Value = 7057.857
$ws.Range("I122").Value = 8901.4
$ws.Range("J122").Value = 2449
$ws.Range("K122").Value = 26704.2
$ws.Range("L122").Value = 7347
$ws.Range("M122").Value = -24254.2
$ws.Range("N122").Value = -12247

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 3500.5789
$ws.Range("I16").Value = 3758.1714
$ws.Range("J16").Value = 495.33334
$ws.Range("K16").Value = 3758.1714
$ws.Range("L16").Value = 495.33334
$ws.Range("M16").Value = -3588.1714
$ws.Range("N16").Value = -835.33334
# Row 32
$ws.Range("H32").Value = 14253.25
$ws.Range("J32").Value = 40000
$ws.Range("L32").Value = 40000
$ws.Range("N32").Value = -40634
# Row 82
$ws.Range("H82").Value = 1925.625
$ws.Range("I82").Value = 1697.25
$ws.Range("J82").Value = 2154
$ws.Range("K82").Value = 1697.25
$ws.Range("L82").Value = 2154
$ws.Range("M82").Value = -1336.25
$ws.Range("N82").Value = -2876
# Row 85
$ws.Range("H85").Value = 1925.625
$ws.Range("I85").Value = 1697.25
$ws.Range("J85").Value = 2154
$ws.Range("K85").Value = 1697.25
$ws.Range("L85").Value = 2154
$ws.Range("M85").Value = -449.25
$ws.Range("N85").Value = -4650
# Row 93
$ws.Range("H93").Value = 1248.25
$ws.Range("I93").Value = 997.6667
$ws.Range("J93").Value = 2000
$ws.Range("K93").Value = 997.6667
$ws.Range("L93").Value = 2000
$ws.Range("M93").Value = 250.3333
$ws.Range("N93").Value = -4496
# Row 133
$ws.Range("H133").Value = 95989
$ws.Range("J133").Value = 95989
$ws.Range("L133").Value = 95989
$ws.Range("N133").Value = -101049
# Row 136
$ws.Range("H136").Value = 29445.111
$ws.Range("I136").Value = 3305.8
$ws.Range("J136").Value = 160141.67
$ws.Range("K136").Value = 9917.400000000001
$ws.Range("L136").Value = 480425.01
$ws.Range("M136").Value = -7367.400000000001
$ws.Range("N136").Value = -485525.01

$ws = $wb.Worksheets.Item("WVR")
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
# Row 81
$ws.Range("H81").Value = 38340.668
$ws.Range("J81").Value = 1989.2
$ws.Range("L81").Value = 3978.4
$ws.Range("N81").Value = -6100.4
# Row 84
$ws.Range("H84").Value = 38340.668
$ws.Range("J84").Value = 1989.2
$ws.Range("L84").Value = 19892
$ws.Range("N84").Value = -30500

